# Auto-generated script to update market-price derived columns (H-N)
# across multiple worksheets, per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 835.55554
$ws.Range("I33").Value = 420
$ws.Range("K33").Value = 420
$ws.Range("M33").Value = -191
$ws.Range("H62").Value = 3664.3333
$ws.Range("I62").Value = 3139.8572
$ws.Range("K62").Value = 3139.8572
$ws.Range("M62").Value = -2515.8572
$ws.Range("H65").Value = 3664.3333
$ws.Range("I65").Value = 3139.8572
$ws.Range("K65").Value = 15699.286
$ws.Range("M65").Value = -12579.286
$ws.Range("H86").Value = 3374.75
$ws.Range("I86").Value = 2441.6155
$ws.Range("K86").Value = 2441.6155
$ws.Range("M86").Value = -1318.6155
$ws.Range("H89").Value = 3374.75
$ws.Range("I89").Value = 2441.6155
$ws.Range("K89").Value = 12208.0775
$ws.Range("M89").Value = -6592.077499999999
$ws.Range("H132").Value = 21313.162
$ws.Range("I132").Value = 3725.9048
$ws.Range("J132").Value = 58246.4
$ws.Range("K132").Value = 11177.7144
$ws.Range("L132").Value = 174739.2
$ws.Range("M132").Value = -8647.714399999999
$ws.Range("N132").Value = -179799.2
$ws.Range("H138").Value = 3659.125
$ws.Range("I138").Value = 2022.4286
$ws.Range("J138").Value = 4333.0586
$ws.Range("K138").Value = 6067.2858
$ws.Range("L138").Value = 12999.1758
$ws.Range("M138").Value = -927.2857999999997
$ws.Range("N138").Value = -23279.1758

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12992587
$ws.Range("I32").Value = 13519071
$ws.Range("K32").Value = 13519071
$ws.Range("M32").Value = -13518784
$ws.Range("H61").Value = 2272.6843
$ws.Range("I61").Value = 2065.8
$ws.Range("K61").Value = 2065.8
$ws.Range("M61").Value = -1853.8
$ws.Range("H63").Value = 9324.875
$ws.Range("J63").Value = 13866.667
$ws.Range("L63").Value = 13866.667
$ws.Range("N63").Value = -15238.667
$ws.Range("H66").Value = 9324.875
$ws.Range("J66").Value = 13866.667
$ws.Range("L66").Value = 69333.33499999999
$ws.Range("N66").Value = -76197.33499999999
$ws.Range("H97").Value = 2413.4375
$ws.Range("I97").Value = 1186
$ws.Range("K97").Value = 1186
$ws.Range("M97").Value = -690
$ws.Range("H102").Value = 1586.6
$ws.Range("I102").Value = 1651.8889
$ws.Range("J102").Value = 999
$ws.Range("K102").Value = 1651.8889
$ws.Range("L102").Value = 999
$ws.Range("M102").Value = -29.88889999999992
$ws.Range("N102").Value = -4243
$ws.Range("H122").Value = 2397.1428
$ws.Range("I122").Value = 2130
$ws.Range("K122").Value = 6390
$ws.Range("M122").Value = -3940
$ws.Range("H136").Value = 2272.6843
$ws.Range("I136").Value = 2065.8
$ws.Range("K136").Value = 6197.400000000001
$ws.Range("M136").Value = -3647.400000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3447.6333
$ws.Range("I86").Value = 2522.8572
$ws.Range("J86").Value = 5605.4443
$ws.Range("K86").Value = 2522.8572
$ws.Range("L86").Value = 5605.4443
$ws.Range("M86").Value = -1399.8572
$ws.Range("N86").Value = -7851.4443
$ws.Range("H89").Value = 3447.6333
$ws.Range("I89").Value = 2522.8572
$ws.Range("J89").Value = 5605.4443
$ws.Range("K89").Value = 12614.286
$ws.Range("L89").Value = 28027.2215
$ws.Range("M89").Value = -6998.286
$ws.Range("N89").Value = -39259.2215
$ws.Range("H94").Value = 3957.5386
$ws.Range("I94").Value = 3957.5386
$ws.Range("K94").Value = 3957.5386
$ws.Range("M94").Value = -3506.5386

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 12548.143
$ws.Range("J16").Value = 7247.6665
$ws.Range("L16").Value = 7247.6665
$ws.Range("N16").Value = -7821.6665
$ws.Range("H36").Value = 12730
$ws.Range("I36").Value = 12730
$ws.Range("K36").Value = 12730
$ws.Range("M36").Value = -12342
$ws.Range("H40").Value = 12730
$ws.Range("I40").Value = 12730
$ws.Range("K40").Value = 12730
$ws.Range("M40").Value = -12570
$ws.Range("H86").Value = 39755.332
$ws.Range("I86").Value = 80750.5
$ws.Range("J86").Value = 19257.75
$ws.Range("K86").Value = 80750.5
$ws.Range("L86").Value = 19257.75
$ws.Range("M86").Value = -79627.5
$ws.Range("N86").Value = -21503.75
$ws.Range("H89").Value = 39755.332
$ws.Range("I89").Value = 80750.5
$ws.Range("J89").Value = 19257.75
$ws.Range("K89").Value = 403752.5
$ws.Range("L89").Value = 96288.75
$ws.Range("M89").Value = -398136.5
$ws.Range("N89").Value = -107520.75
$ws.Range("H113").Value = 12548.143
$ws.Range("J113").Value = 7247.6665
$ws.Range("L113").Value = 7247.6665
$ws.Range("N113").Value = -11587.6665

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 1099.2
$ws.Range("J92").Value = 1274.25
$ws.Range("L92").Value = 3822.75
$ws.Range("N92").Value = -6318.75
$ws.Range("H123").Value = 710
$ws.Range("I123").Value = 710
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 2130
$ws.Range("L123").ClearContents()
$ws.Range("N123").Value = 0
$ws.Range("M123").Value = 320
$ws.Range("H131").Value = 3569.2354
$ws.Range("I131").Value = 2728.2856
$ws.Range("J131").Value = 3787.2593
$ws.Range("K131").Value = 8184.8568
$ws.Range("L131").Value = 11361.7779
$ws.Range("M131").Value = -3144.8568
$ws.Range("N131").Value = -21441.7779

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 205.25
$ws.Range("I2").Value = 130.75
$ws.Range("K2").Value = 130.75
$ws.Range("M2").Value = -17.75
$ws.Range("H70").Value = 189351.67
$ws.Range("J70").Value = 4999.75
$ws.Range("L70").Value = 4999.75
$ws.Range("N70").Value = -5539.75
$ws.Range("H73").Value = 189351.67
$ws.Range("J73").Value = 4999.75
$ws.Range("L73").Value = 4999.75
$ws.Range("N73").Value = -6871.75
$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -33744
$ws.Range("H102").Value = 6413.125
$ws.Range("I102").Value = 2721.6
$ws.Range("K102").Value = 2721.6
$ws.Range("M102").Value = -1099.6
$ws.Range("H122").Value = 9199.200000000001
$ws.Range("I122").Value = 9332
$ws.Range("K122").Value = 27996
$ws.Range("M122").Value = -25546
$ws.Range("H126").Value = 5886
$ws.Range("J126").Value = 6894
$ws.Range("L126").Value = 20682
$ws.Range("N126").Value = -25622
$ws.Range("H132").Value = 1898.6842
$ws.Range("I132").Value = 1605
$ws.Range("K132").Value = 4815
$ws.Range("M132").Value = -2285
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").ClearContents()
$ws.Range("N138").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7319.85
$ws.Range("I7").Value = 8178.9165
$ws.Range("J7").Value = 6031.25
$ws.Range("K7").Value = 8178.9165
$ws.Range("L7").Value = 6031.25
$ws.Range("M7").Value = -8066.9165
$ws.Range("N7").Value = -6255.25
$ws.Range("I16").Value = 1992.6666
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 1992.6666
$ws.Range("L16").Value = 500
$ws.Range("M16").Value = -1822.6666
$ws.Range("N16").Value = -840
$ws.Range("H40").Value = 13216.6
$ws.Range("I40").Value = 15251.375
$ws.Range("K40").Value = 15251.375
$ws.Range("M40").Value = -15115.375
$ws.Range("H55").Value = 3577
$ws.Range("I55").Value = 2983.625
$ws.Range("J55").Value = 4368.1665
$ws.Range("K55").Value = 2983.625
$ws.Range("L55").Value = 4368.1665
$ws.Range("M55").Value = -2810.625
$ws.Range("N55").Value = -4714.1665
$ws.Range("H68").Value = 3503.5
$ws.Range("I68").Value = 3306
$ws.Range("J68").Value = 3780
$ws.Range("K68").Value = 3306
$ws.Range("L68").Value = 3780
$ws.Range("M68").Value = -2557
$ws.Range("N68").Value = -5278
$ws.Range("H71").Value = 3503.5
$ws.Range("I71").Value = 3306
$ws.Range("J71").Value = 3780
$ws.Range("K71").Value = 16530
$ws.Range("L71").Value = 18900
$ws.Range("M71").Value = -12786
$ws.Range("N71").Value = -26388
$ws.Range("H116").Value = 249992.25
$ws.Range("J116").Value = 249992.25
$ws.Range("L116").Value = 249992.25
$ws.Range("N116").Value = -259170.25
$ws.Range("H122").Value = 4537.838
$ws.Range("I122").Value = 1983.0555
$ws.Range("J122").Value = 6958.1577
$ws.Range("K122").Value = 5949.166499999999
$ws.Range("L122").Value = 20874.4731
$ws.Range("M122").Value = -3499.166499999999
$ws.Range("N122").Value = -25774.4731
$ws.Range("H126").Value = 7319.85
$ws.Range("I126").Value = 8178.9165
$ws.Range("J126").Value = 6031.25
$ws.Range("K126").Value = 24536.7495
$ws.Range("L126").Value = 18093.75
$ws.Range("M126").Value = -22066.7495
$ws.Range("N126").Value = -23033.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2581.0833
$ws.Range("I122").Value = 2581.0833
$ws.Range("K122").Value = 7743.249899999999
$ws.Range("M122").Value = -5293.249899999999
$ws.Range("H126").Value = 2826
$ws.Range("I126").Value = 2117.1667
$ws.Range("J126").Value = 4952.5
$ws.Range("K126").Value = 6351.500100000001
$ws.Range("L126").Value = 14857.5
$ws.Range("M126").Value = -3881.500100000001
$ws.Range("N126").Value = -19797.5
$ws.Range("H132").Value = 1522.0209
$ws.Range("I132").Value = 1536.4318
$ws.Range("K132").Value = 4609.2954
$ws.Range("M132").Value = -2079.2954
$ws.Range("H136").Value = 1177.6383
$ws.Range("I136").Value = 654.875
$ws.Range("J136").Value = 4164.857
$ws.Range("K136").Value = 1964.625
$ws.Range("L136").Value = 12494.571
$ws.Range("M136").Value = 585.375
$ws.Range("N136").Value = -17594.571
